$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new trailing columns: Wins / Losses / Ties ------------------
# Column AC (29) already holds the last existing header ("Unnamed: 28").
# Copy its header formatting (bold, centered, bordered) onto the three new
# header cells AD1:AF1 before writing their labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Populate the season record for every player row (2-57) ---------------
# Every row in this sheet shares the same team season record: 55 wins,
# 107 losses, 0 ties.
$ws.Range("AD2:AD57").Value = 55
$ws.Range("AE2:AE57").Value = 107
$ws.Range("AF2:AF57").Value = 0
